$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing a Text number format so
# numeric-looking strings (and percent-looking strings) are preserved exactly
# as text, matching how the source data is stored (inline/shared strings).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2
Set-TextValue "D2" "301.62"
Set-TextValue "E2" "-0.58%"

# Row 3
Set-TextValue "E3" "-2.03%"

# Row 4
Set-TextValue "D4" "5.091"
Set-TextValue "E4" "-3.25%"

# Row 5
Set-TextValue "D5" "0.07362"
Set-TextValue "E5" "-2.26%"

# Row 6
Set-TextValue "D6" "2.307"
Set-TextValue "E6" "53.68%"

# Row 7
Set-TextValue "D7" "7.945"
Set-TextValue "E7" "0.43%"

# Row 8
Set-TextValue "D8" "3.783"
Set-TextValue "E8" "-0.81%"

# Row 9
Set-TextValue "D9" "0.9187"
Set-TextValue "E9" "-0.24%"

# Row 10
Set-TextValue "D10" "0.1705"
Set-TextValue "E10" "0.61%"

# Row 11
Set-TextValue "D11" "0.07493"
Set-TextValue "E11" "-6.40%"

# Row 12
Set-TextValue "D12" "0.08170"
Set-TextValue "E12" "1.03%"

# Row 13
Set-TextValue "D13" "0.03025"
Set-TextValue "E13" "-0.36%"

# Row 14
Set-TextValue "D14" "0.09940"
Set-TextValue "E14" "0.27%"

# Row 15
Set-TextValue "D15" "0.001494"
Set-TextValue "E15" "-1.75%"

# Row 16
Set-TextValue "D16" "0.006084"
Set-TextValue "E16" "-6.94%"

# Row 17
Set-TextValue "E17" "-0.06%"

# Row 18
Set-TextValue "D18" "2.222"
Set-TextValue "E18" "-0.24%"

# Row 20
Set-TextValue "D20" "0.1346"
Set-TextValue "E20" "2.97%"

# Row 21
Set-TextValue "D21" "4.647"
Set-TextValue "E21" "3.29%"

# Row 22
Set-TextValue "D22" "0.04646"
Set-TextValue "E22" "0.92%"

# Row 23
Set-TextValue "D23" "0.1566"
Set-TextValue "E23" "-3.22%"

# Row 24
Set-TextValue "D24" "0.001226"
Set-TextValue "E24" "0.80%"

# Row 25
Set-TextValue "D25" "0.004491"
Set-TextValue "E25" "0.78%"

# Row 26
Set-TextValue "D26" "0.0001299"
Set-TextValue "E26" "-7.20%"

# Row 27
Set-TextValue "E27" "50.45%"

# Row 39
Set-TextValue "D39" "0.01745"
Set-TextValue "E39" "2.02%"

# Row 40
Set-TextValue "D40" "0.04526"
Set-TextValue "E40" "0.83%"

# Row 41
Set-TextValue "D41" "0.007201"
Set-TextValue "E41" "3.96%"

# Row 42
Set-TextValue "D42" "0.1348"
Set-TextValue "E42" "-0.18%"

# Row 43
Set-TextValue "D43" "0.002228"
Set-TextValue "E43" "4.14%"

# Row 44
Set-TextValue "D44" "0.01074"
Set-TextValue "E44" "-16.27%"

# Row 45
Set-TextValue "D45" "0.00006293"
Set-TextValue "E45" "2.07%"

# Row 46/47: BOLO and CoinbaseStockToken swap positions, with refreshed data
Set-TextValue "B46" "CoinbaseStockToken"
Set-TextValue "C46" "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue "D46" "0.009992"
Set-TextValue "E46" "-33.28%"

Set-TextValue "B47" "BOLO"
Set-TextValue "C47" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue "D47" "0.8085"
Set-TextValue "E47" "13.55%"
